$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to match new, longer explanatory text (approx. target 67.29 chars;
# engine quantizes ColumnWidth to 1/6-character steps, so 66.5 is the closest achievable input).
$ws.Columns("A").ColumnWidth = 66.5

# Reuse the existing highlighted-label style (fill, from column A row 1) for the new label/value rows
# by copying formats only (keeps the same cellXf/fill index instead of creating a new style slot).
$ws.Range("A1:B1").Copy()

$ws.Range("A17:B17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Thuận mùa sinh (Mệnh vượng vào mùa Xuân)"
$ws.Range("B17").Value = "Tăng độ số may mắn."

$ws.Range("A18:B18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Thuận mùa sinh (Mệnh vượng vào mùa Hạ)"
$ws.Range("B18").Value = "Tăng độ số may mắn."

$ws.Range("A19:B19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Thuận mùa sinh (Mệnh vượng vào mùa Thu)"
$ws.Range("B19").Value = "Tăng độ số may mắn."

$ws.Range("A20:B20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Thuận mùa sinh (Mệnh vượng vào mùa Đông)"
$ws.Range("B20").Value = "Tăng độ số may mắn."

$ws.Range("A21:B21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Không thuận mùa sinh (Mệnh không vượng/tướng vào mùa Xuân)"
$ws.Range("B21").Value = "Giảm độ số may mắn."

$ws.Range("A22:B22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Không thuận mùa sinh (Mệnh không vượng/tướng vào mùa Hạ)"
$ws.Range("B22").Value = "Giảm độ số may mắn."

$ws.Range("A23:B23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Không thuận mùa sinh (Mệnh không vượng/tướng vào mùa Thu)"
$ws.Range("B23").Value = "Giảm độ số may mắn."

$ws.Range("A24:B24").PasteSpecial(-4122)
$ws.Range("A24").Value = "Không thuận mùa sinh (Mệnh không vượng/tướng vào mùa Đông)"
$ws.Range("B24").Value = "Giảm độ số may mắn."

$ws.Range("A25:B25").PasteSpecial(-4122)
$ws.Range("A25").Value = "Thuận giờ sinh (Mệnh vượng giờ Dậu mùa Xuân)"
$ws.Range("B25").Value = "Tăng độ số may mắn."

$ws.Range("A26:B26").PasteSpecial(-4122)
$ws.Range("A26").Value = "Thuận giờ sinh (Mệnh vượng giờ Thìn mùa Xuân)"
$ws.Range("B26").Value = "Tăng độ số may mắn."

$ws.Range("A27:B27").PasteSpecial(-4122)
$ws.Range("A27").Value = "Thuận giờ sinh (Mệnh vượng giờ Tỵ mùa Xuân)"
$ws.Range("B27").Value = "Tăng độ số may mắn."

$ws.Range("A28:B28").PasteSpecial(-4122)
$ws.Range("A28").Value = "Thuận giờ sinh (Mệnh vượng giờ Mão mùa Hạ)"
$ws.Range("B28").Value = "Tăng độ số may mắn."

$ws.Range("A29:B29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Thuận giờ sinh (Mệnh vượng giờ Mùi mùa Hạ)"
$ws.Range("B29").Value = "Tăng độ số may mắn."

$ws.Range("A30:B30").PasteSpecial(-4122)
$ws.Range("A30").Value = "Thuận giờ sinh (Mệnh vượng giờ Hợi mùa Hạ)"
$ws.Range("B30").Value = "Tăng độ số may mắn."

$ws.Range("A31:B31").PasteSpecial(-4122)
$ws.Range("A31").Value = "Thuận giờ sinh (Mệnh vượng giờ Ngọ mùa Thu)"
$ws.Range("B31").Value = "Tăng độ số may mắn."

$ws.Range("A32:B32").PasteSpecial(-4122)
$ws.Range("A32").Value = "Thuận giờ sinh (Mệnh vượng giờ Thân mùa Thu)"
$ws.Range("B32").Value = "Tăng độ số may mắn."

$ws.Range("A33:B33").PasteSpecial(-4122)
$ws.Range("A33").Value = "Thuận giờ sinh (Mệnh vượng giờ Dần mùa Thu)"
$ws.Range("B33").Value = "Tăng độ số may mắn."

$ws.Range("A34:B34").PasteSpecial(-4122)
$ws.Range("A34").Value = "Thuận giờ sinh (Mệnh vượng giờ Sửu mùa Đông)"
$ws.Range("B34").Value = "Tăng độ số may mắn."

$ws.Range("A35:B35").PasteSpecial(-4122)
$ws.Range("A35").Value = "Thuận giờ sinh (Mệnh vượng giờ Tuất mùa Đông)"
$ws.Range("B35").Value = "Tăng độ số may mắn."

$ws.Range("A36:B36").PasteSpecial(-4122)
$ws.Range("A36").Value = "Thuận giờ sinh (Mệnh vượng giờ Tý mùa Đông)"
$ws.Range("B36").Value = "Tăng độ số may mắn."

$ws.Range("A37:B37").PasteSpecial(-4122)
$ws.Range("A37").Value = "Không thuận giờ sinh theo mùa sinh"
$ws.Range("B37").Value = "Giảm độ số may mắn."

$ws.Range("A38:B38").PasteSpecial(-4122)
$ws.Range("A38").Value = "Năm sinh và giờ sinh phạm tối độc"
$ws.Range("B38").Value = "Năm sinh và giờ sinh phạm tối độc: số cách biệt cha mẹ, khó ở lâu dài với cha mẹ."

$ws.Range("A39:B39").PasteSpecial(-4122)
$ws.Range("A39").Value = "Năm sinh và giờ sinh phạm hình khắc cha lúc nhỏ tuổi"
$ws.Range("B39").Value = "Năm sinh và giờ sinh phạm hình khắc cha lúc nhỏ tuổi. Nếu qua 16 tuổi mà cha vẫn còn thì hình khắc đã tiêu, cha con có thể chung sống lâu dài."

$ws.Range("A40:B40").PasteSpecial(-4122)
$ws.Range("A40").Value = "Năm sinh và giờ sinh không phạm tối độc"
$ws.Range("B40").Value = "Bình thường, không hình thương khắc hại cha mẹ."

$ws.Range("A41:B41").PasteSpecial(-4122)
$ws.Range("A41").Value = "Năm sinh và giờ sinh phạm hình khắc mẹ"
$ws.Range("B41").Value = "Mẹ có thể chết trước cha."

# Row 42: a trailing styled-but-empty label cell (no value), same style as the rest of column A.
$ws.Range("A1").Copy()
$ws.Range("A42").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Selection / scroll position, matching the saved view in the source file.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M42").Select()

# Page orientation, as set in Page Setup.
$ws.PageSetup.Orientation = 1

